$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add participant 4 study data: name in A5 (B5 already has ID 4)
$ws.Range("A5").Value = "Grey Nclayghlin"

# Update the active cell selection to A5 to match the author's last edit position
$ws.Range("A5").Select()
